# Weekly update: two new Chirimoya price observations arrived, so they are
# inserted at the top of the existing data block (rows 16:17), pushing all
# the subsequent rows down by two (old row 16 -> new row 18, ... old row
# 37 -> new row 39). This matches the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 16, shifting rows 16:37 down
# to 18:39 (and carrying the date-formatted style of column D along).
$ws.Rows("16:17").Insert()

# --- New row 16 ---
$ws.Cells.Item(16, 1).Value2 = 5
$ws.Cells.Item(16, 2).Value2 = 'Macroferia Regional de Talca'
$ws.Cells.Item(16, 3).Value2 = 'Maule'
$ws.Cells.Item(16, 4).Value2 = 44484
$ws.Cells.Item(16, 5).Value2 = 7
$ws.Cells.Item(16, 6).Value2 = 'Fruta'
$ws.Cells.Item(16, 7).Value2 = 100107
$ws.Cells.Item(16, 8).Value2 = 'Otros'
$ws.Cells.Item(16, 9).Value2 = 100107002
$ws.Cells.Item(16, 10).Value2 = 'Chirimoya'
$ws.Cells.Item(16, 11).Value2 = 'Cultivar IV Región'
$ws.Cells.Item(16, 12).Value2 = 'Primera'
$ws.Cells.Item(16, 13).Value2 = 120
$ws.Cells.Item(16, 14).Value2 = 25000
$ws.Cells.Item(16, 15).Value2 = 25000
$ws.Cells.Item(16, 16).Value2 = 25000
$ws.Cells.Item(16, 17).Value2 = '$/bandeja 10 kilos'
$ws.Cells.Item(16, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(16, 19).Value2 = 2500
$ws.Cells.Item(16, 20).Value2 = 10

# --- New row 17 ---
$ws.Cells.Item(17, 1).Value2 = 5
$ws.Cells.Item(17, 2).Value2 = 'Macroferia Regional de Talca'
$ws.Cells.Item(17, 3).Value2 = 'Maule'
$ws.Cells.Item(17, 4).Value2 = 44484
$ws.Cells.Item(17, 5).Value2 = 7
$ws.Cells.Item(17, 6).Value2 = 'Fruta'
$ws.Cells.Item(17, 7).Value2 = 100107
$ws.Cells.Item(17, 8).Value2 = 'Otros'
$ws.Cells.Item(17, 9).Value2 = 100107002
$ws.Cells.Item(17, 10).Value2 = 'Chirimoya'
$ws.Cells.Item(17, 11).Value2 = 'Cultivar IV Región'
$ws.Cells.Item(17, 12).Value2 = 'Segunda'
$ws.Cells.Item(17, 13).Value2 = 100
$ws.Cells.Item(17, 14).Value2 = 22000
$ws.Cells.Item(17, 15).Value2 = 22000
$ws.Cells.Item(17, 16).Value2 = 22000
$ws.Cells.Item(17, 17).Value2 = '$/bandeja 10 kilos'
$ws.Cells.Item(17, 18).Value2 = 'Provincia de Limarí'
$ws.Cells.Item(17, 19).Value2 = 2200
$ws.Cells.Item(17, 20).Value2 = 10
